$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to store its Value as literal text (not auto-converted to a number),
# while leaving the cell style index unchanged (matches original un-styled cells).
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "37.870.28"
Set-TextValue $ws.Range("E2") "  -0.68%  "
Set-TextValue $ws.Range("D3") "2.049.10"
Set-TextValue $ws.Range("E3") "  -0.21%  "
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "228.64"
Set-TextValue $ws.Range("E5") "  -0.25%  "
Set-TextValue $ws.Range("D6") "0.609"
Set-TextValue $ws.Range("E6") "  -1.13%  "
Set-TextValue $ws.Range("D7") "61.00"
Set-TextValue $ws.Range("E7") "  +0.56%  "
Set-TextValue $ws.Range("E8") "  -0.01%  "
Set-TextValue $ws.Range("D9") "0.377"
Set-TextValue $ws.Range("E9") "  -1.83%  "
Set-TextValue $ws.Range("D10") "0.0821"
Set-TextValue $ws.Range("E10") "  -0.41%  "
Set-TextValue $ws.Range("D11") "0.105"
Set-TextValue $ws.Range("E11") "  +0.50%  "
Set-TextValue $ws.Range("D12") "14.78"
Set-TextValue $ws.Range("E12") "  -0.03%  "
Set-TextValue $ws.Range("D13") "2.355.26"
Set-TextValue $ws.Range("E13") "  -0.12%  "
Set-TextValue $ws.Range("D14") "21.12"
Set-TextValue $ws.Range("E14") "  -0.14%  "
Set-TextValue $ws.Range("D15") "0.778"
Set-TextValue $ws.Range("E15") "  +2.77%  "
Set-TextValue $ws.Range("E16") "  -2.55%  "
Set-TextValue $ws.Range("D17") "2.020.05"
Set-TextValue $ws.Range("E17") "  -1.38%  "
Set-TextValue $ws.Range("D18") "37.861.40"
Set-TextValue $ws.Range("E18") "  -0.60%  "
Set-TextValue $ws.Range("D19") "69.65"
Set-TextValue $ws.Range("E19") "  -0.18%  "
Set-TextValue $ws.Range("D20") "5.92"
Set-TextValue $ws.Range("E20") "  -5.31%  "
Set-TextValue $ws.Range("D21") "0.0₃0824"
Set-TextValue $ws.Range("E21") "  -1.35%  "
Set-TextValue $ws.Range("D22") "223.95"
Set-TextValue $ws.Range("E22") "  -0.67%  "
Set-TextValue $ws.Range("D23") "1.00"
Set-TextValue $ws.Range("E23") "  +0.06%  "
Set-TextValue $ws.Range("D24") "2.40"
Set-TextValue $ws.Range("E24") "  -1.28%  "
Set-TextValue $ws.Range("E25") "  +2.57%  "
Set-TextValue $ws.Range("B26") "Monero"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "168.23"
Set-TextValue $ws.Range("E26") "  +1.53%  "
Set-TextValue $ws.Range("B27") "Cosmos"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "9.38"
Set-TextValue $ws.Range("E27") "  +1.64%  "
Set-TextValue $ws.Range("D28") "0.130"
Set-TextValue $ws.Range("E28") "  -2.66%  "
Set-TextValue $ws.Range("D29") "18.85"
Set-TextValue $ws.Range("E29") "  -0.43%  "
Set-TextValue $ws.Range("E30") "  -1.53%  "
Set-TextValue $ws.Range("D31") "0.120"
Set-TextValue $ws.Range("E31") "  +0.19%  "
Set-TextValue $ws.Range("D32") "2.22"
Set-TextValue $ws.Range("E32") "  +8.13%  "
Set-TextValue $ws.Range("E33") "  -2.60%  "
Set-TextValue $ws.Range("E34") "  -0.92%  "
Set-TextValue $ws.Range("D35") "0.0606"
Set-TextValue $ws.Range("E35") "  +0.38%  "
Set-TextValue $ws.Range("D36") "6.55"
Set-TextValue $ws.Range("E36") "  +2.42%  "
Set-TextValue $ws.Range("D37") "2.35"
Set-TextValue $ws.Range("E37") "  +3.31%  "
Set-TextValue $ws.Range("D38") "3.41"
Set-TextValue $ws.Range("E38") "  +4.35%  "
Set-TextValue $ws.Range("E39") "  -0.05%  "
Set-TextValue $ws.Range("D40") "18.01"
Set-TextValue $ws.Range("E40") "  +6.08%  "
Set-TextValue $ws.Range("D41") "1.540.03"
Set-TextValue $ws.Range("E41") "  +1.38%  "
Set-TextValue $ws.Range("E42") "  +1.04%  "
Set-TextValue $ws.Range("D43") "96.36"
Set-TextValue $ws.Range("E43") "  -0.84%  "
Set-TextValue $ws.Range("E44") "  -1.66%  "
Set-TextValue $ws.Range("D45") "0.0912"
Set-TextValue $ws.Range("E45") "  -1.71%  "
Set-TextValue $ws.Range("E46") "  -1.73%  "
Set-TextValue $ws.Range("D47") "4.09"
Set-TextValue $ws.Range("E47") "  +2.31%  "
Set-TextValue $ws.Range("E48") "  +0.00%  "
Set-TextValue $ws.Range("E49") "  -0.06%  "
Set-TextValue $ws.Range("D50") "7.05"
Set-TextValue $ws.Range("E50") "  +0.30%  "
Set-TextValue $ws.Range("D51") "2.243.61"
Set-TextValue $ws.Range("E51") "  -0.10%  "
